# Commit: "new task defaults status and category to the task above it"
#
# The todo item with Id=8 ("when new task is created, default set its
# drop downs to the same values as the one above it") has been
# implemented, so it moves from the "Active" sheet to the "Inactive"
# sheet: Status becomes "Done" and a completion date ("3/3/2018") is
# recorded in the "Done" column. Category/Created stay the same as
# before. It is inserted as the first data row (row 2, right after the
# header) of "Inactive", and removed from its old spot (row 3) on
# "Active" - every other row shifts accordingly.

$wb = $excel.ActiveWorkbook
$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")

# Make room for the newly-completed task at the top of the Inactive list
# and reset the inserted row's formatting back to the plain data-row
# style (Insert() otherwise inherits the bold header style from row 1).
$wsInactive.Rows.Item(2).Insert()
$wsInactive.Range("A2:F2").ClearFormats()

$wsInactive.Range("A2").Value = 8
$wsInactive.Range("B2").Value = "when new task is created, default set its drop downs to the same values as the one above it"
$wsInactive.Range("C2").Value = "Done"
$wsInactive.Range("D2").Value = "Feature"

# The Created/Done columns store dates as plain text (e.g. "12/1/2017"),
# not real date serials. Typing a date-shaped string straight into
# .Value would get auto-converted to a date, so instead copy the value
# + number format from existing text cells that already hold the exact
# same string.
$wsActive.Range("E2").Copy()
$wsInactive.Range("E2").PasteSpecial(-4163)   # xlPasteValuesAndNumberFormats: "12/1/2017"

$wsInactive.Range("F3").Copy()
$wsInactive.Range("F2").PasteSpecial(-4163)   # xlPasteValuesAndNumberFormats: "3/3/2018"

# Remove the now-completed task from the Active sheet (it was row 3,
# Id=8); remaining rows shift up.
$wsActive.Rows.Item(3).Delete()
